$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2958.3
$ws.Range("I32").Value = 2422.5
$ws.Range("K32").Value = 2422.5
$ws.Range("M32").Value = -2096.5
$ws.Range("H51").Value = 7708
$ws.Range("I51").Value = 9115.571
$ws.Range("K51").Value = 9115.571
$ws.Range("M51").Value = -8631.571
$ws.Range("H55").Value = 758.5454999999999
$ws.Range("I55").Value = 590
$ws.Range("K55").Value = 590
$ws.Range("M55").Value = -376
$ws.Range("H62").Value = 11890.85
$ws.Range("J62").Value = 9524
$ws.Range("L62").Value = 9524
$ws.Range("N62").Value = -10772
$ws.Range("H65").Value = 11890.85
$ws.Range("J65").Value = 9524
$ws.Range("L65").Value = 47620
$ws.Range("N65").Value = -53860
$ws.Range("H76").Value = 4418.1665
$ws.Range("I76").Value = 4119.25
$ws.Range("K76").Value = 4119.25
$ws.Range("M76").Value = -3804.25
$ws.Range("H79").Value = 4418.1665
$ws.Range("I79").Value = 4119.25
$ws.Range("K79").Value = 4119.25
$ws.Range("M79").Value = -3027.25
$ws.Range("H128").Value = 110000
$ws.Range("J128").Value = 110000
$ws.Range("L128").Value = 110000
$ws.Range("N128").Value = -119960
$ws.Range("H132").Value = 10425.475
$ws.Range("I132").Value = 9420.361000000001
$ws.Range("J132").Value = 19471.5
$ws.Range("K132").Value = 28261.083
$ws.Range("L132").Value = 58414.5
$ws.Range("M132").Value = -25731.083
$ws.Range("N132").Value = -63474.5
$ws.Range("H137").Value = 1233.5834
$ws.Range("I137").Value = 740.4
$ws.Range("J137").Value = 3699.5
$ws.Range("K137").Value = 2221.2
$ws.Range("L137").Value = 11098.5
$ws.Range("M137").Value = 328.8000000000002
$ws.Range("N137").Value = -16198.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13976.452
$ws.Range("I32").Value = 14011.728
$ws.Range("J32").Value = 13533
$ws.Range("K32").Value = 14011.728
$ws.Range("L32").Value = 13533
$ws.Range("M32").Value = -13724.728
$ws.Range("N32").Value = -14107
$ws.Range("H61").Value = 1664.5
$ws.Range("I61").Value = 1664.5
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1664.5
$ws.Range("L61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("N61").ClearContents()
$ws.Range("H136").Value = 1664.5
$ws.Range("I136").Value = 1664.5
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 4993.5
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1303.5
$ws.Range("J22").Value = 2622.8333
$ws.Range("L22").Value = 2622.8333
$ws.Range("N22").Value = -3322.8333
$ws.Range("H99").Value = 2379.2144
$ws.Range("J99").Value = 2845.2856
$ws.Range("L99").Value = 2845.2856
$ws.Range("N99").Value = -5841.2856
$ws.Range("H107").Value = 479.73685
$ws.Range("I107").Value = 462.64706
$ws.Range("J107").Value = 625
$ws.Range("K107").Value = 462.64706
$ws.Range("L107").Value = 625
$ws.Range("M107").Value = 1457.35294
$ws.Range("N107").Value = -4465
$ws.Range("H126").Value = 2379.2144
$ws.Range("J126").Value = 2845.2856
$ws.Range("L126").Value = 8535.856800000001
$ws.Range("N126").Value = -13475.8568

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H130").Value = 7595.1816
$ws.Range("I130").Value = 7068.75
$ws.Range("K130").Value = 21206.25
$ws.Range("M130").Value = -16186.25
$ws.Range("H131").Value = 5016674.5
$ws.Range("I131").Value = 12962.667
$ws.Range("K131").Value = 38888.001
$ws.Range("M131").Value = -33848.001
$ws.Range("H137").Value = 4154.4546
$ws.Range("J137").Value = 5125
$ws.Range("L137").Value = 15375
$ws.Range("N137").Value = -25575

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 92.666664
$ws.Range("I2").Value = 91.2
$ws.Range("K2").Value = 91.2
$ws.Range("M2").Value = 21.8
$ws.Range("H102").Value = 2656.5557
$ws.Range("I102").Value = 929.5
$ws.Range("K102").Value = 929.5
$ws.Range("M102").Value = 692.5
$ws.Range("H132").Value = 55706.58
$ws.Range("I132").Value = 69682.39999999999
$ws.Range("J132").Value = 3297.25
$ws.Range("K132").Value = 209047.2
$ws.Range("L132").Value = 9891.75
$ws.Range("M132").Value = -206517.2
$ws.Range("N132").Value = -14951.75
$ws.Range("H134").Value = 39999.5
$ws.Range("J134").Value = 39999.5
$ws.Range("L134").Value = 119998.5
$ws.Range("N134").Value = -125068.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 535
$ws.Range("I16").Value = 535
$ws.Range("K16").Value = 535
$ws.Range("M16").Value = -365
$ws.Range("H22").Value = 1924.875
$ws.Range("J22").Value = 4899
$ws.Range("L22").Value = 4899
$ws.Range("N22").Value = -5489
$ws.Range("H27").Value = 1924.875
$ws.Range("J27").Value = 4899
$ws.Range("L27").Value = 4899
$ws.Range("N27").Value = -5113
$ws.Range("H46").Value = 7901.476
$ws.Range("I46").Value = 12386
$ws.Range("J46").Value = 5141.769
$ws.Range("K46").Value = 12386
$ws.Range("L46").Value = 5141.769
$ws.Range("M46").Value = -12198
$ws.Range("N46").Value = -5517.769
$ws.Range("H61").Value = 2717.0527
$ws.Range("I61").Value = 2049.9333
$ws.Range("K61").Value = 2049.9333
$ws.Range("M61").Value = -1847.9333
$ws.Range("H113").Value = 2717.0527
$ws.Range("I113").Value = 2049.9333
$ws.Range("K113").Value = 2049.9333
$ws.Range("M113").Value = 120.0666999999999
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").ClearContents()
$ws.Range("N121").ClearContents()
$ws.Range("H122").Value = 3333
$ws.Range("I122").Value = 2735.4482
$ws.Range("K122").Value = 8206.3446
$ws.Range("M122").Value = -5756.3446

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 25054.5
$ws.Range("I11").Value = 3104
$ws.Range("J11").Value = 47005
$ws.Range("K11").Value = 3104
$ws.Range("L11").Value = 47005
$ws.Range("M11").Value = -2962
$ws.Range("N11").Value = -47289
$ws.Range("H108").Value = 55626
$ws.Range("J108").Value = 55626
$ws.Range("L108").Value = 55626
$ws.Range("N108").Value = -63306
$ws.Range("H122").Value = 2743.7827
$ws.Range("J122").Value = 2825.6667
$ws.Range("L122").Value = 8477.000100000001
$ws.Range("N122").Value = -13377.0001
$ws.Range("H136").Value = 4184.409
$ws.Range("I136").Value = 4212.85
$ws.Range("K136").Value = 12638.55
$ws.Range("M136").Value = -10088.55
